$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JMSD200A_23.1.11_18062025")

$startRow = 44

$newRows = @(
    @("JMSD200A_23.1.11", "jiohotstar_25.06.09.3", "subham", "PCM", "src_fmt = SDR", "1920 1080"),
    @("JMSD200A_23.1.11", "jiohotstar_25.06.09.3", "Salaar", "PCM", "src_fmt = SDR", "1920 1080"),
    @("JMSD200A_23.1.11", "jiohotstar_25.06.09.3", "Bhagavanth Kesari", "PCM", "src_fmt = SDR", "1920 1080"),
    @("JMSD200A_23.1.11", "jiohotstar_25.06.09.3", "Captain America Brave new world", "PCM", "src_fmt = invalid", "NA NA")
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($row, $c).Value = $values[$c - 1]
    }
}

$ws.Columns.Item(5).ColumnWidth = 14.1875
